# RPA datasets push 2024-05-30
#
# The "02_38커뮤니케이션(최근일자기준)" sheet gets a new IPO entry
# ("신한글로벌액티브") inserted as row 18. Every row that was at/after the
# old row 18 shifts down by one, and the row that falls off the bottom
# (the former "라메디텍" row, old row 21) is dropped so the table stays
# at 20 data rows (A1:F21 overall, including the header row).
#
# Implemented as a genuine row-insert + fill + trailing-row-delete so the
# existing rows 2-17 (and row 1's header) are left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02_38커뮤니케이션(최근일자기준)")

# Push old rows 18:21 down to 19:22, opening up a blank row 18.
$ws.Rows("18:18").Insert()

# Populate the new row 18 with the 신한글로벌액티브 IPO data.
$ws.Range("A18").Value = "신한글로벌액티브"
$ws.Range("B18").Value = "2024.06.03~06.05"
$ws.Range("C18").Value = "3,000~3,800"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = 70000
$ws.Range("F18").Value = "신한투자증권,한국투자증권"

# The old last row (라메디텍, now shifted to row 22) is dropped so the
# sheet keeps 21 total rows (1 header + 20 data rows).
$ws.Rows("22:22").Delete()
